$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Step 1: locate the "Assignments:" / "None." paragraph pair in the
# Wednesday, March 27, 2024 section (the one right before the Thursday
# heading), and replace the "None." run text with four paragraphs:
#   1. None.  (unchanged text, but with the TOC bookmark now appended)
#   2. a blank ListParagraph line
#   3. "Tasks:" heading line
#   4. the new struck-through task line
# ---------------------------------------------------------------------------

$wedRange = $d.Content
$wedRange.Find.Execute("Wednesday, March 27, 2024") | Out-Null
$wedStart = $wedRange.End

$thuRange = $d.Content
$thuRange.Find.Execute("Thursday, March 28, 2024") | Out-Null
$thuStart = $thuRange.Start

$section = $d.Range($wedStart, $thuStart)
$f = $section.Find
$f.Text = "Assignments:" + [char]13 + "None."
$f.Forward = $true
$f.Wrap = 0
$f.MatchWildcards = $false
$found = $f.Execute()
if (-not $found) {
    throw "Could not find the Wednesday Assignments/None. paragraph"
}

$noneEnd = $section.End
$noneTarget = $d.Range($noneEnd - 5, $noneEnd)

$newXml = "<w:p $wns>" +
            "<w:r><w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr><w:t>None.</w:t></w:r>" +
            "<w:bookmarkStart w:id=`"8`" w:name=`"_Toc162351908`"/>" +
          "</w:p>" +
          "<w:p $wns>" +
            "<w:pPr>" +
              "<w:pStyle w:val=`"ListParagraph`"/>" +
              "<w:keepNext/>" +
              "<w:keepLines/>" +
              "<w:ind w:left=`"1440`"/>" +
              "<w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr>" +
            "</w:pPr>" +
          "</w:p>" +
          "<w:p $wns>" +
            "<w:pPr>" +
              "<w:pStyle w:val=`"ListParagraph`"/>" +
              "<w:keepNext/>" +
              "<w:keepLines/>" +
              "<w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr>" +
              "<w:ind w:left=`"720`"/>" +
              "<w:rPr><w:color w:val=`"E97132`" w:themeColor=`"accent2`"/></w:rPr>" +
            "</w:pPr>" +
            "<w:r>" +
              "<w:rPr><w:color w:val=`"E97132`" w:themeColor=`"accent2`"/><w:u w:val=`"single`"/></w:rPr>" +
              "<w:t>Tasks:</w:t>" +
            "</w:r>" +
          "</w:p>" +
          "<w:p $wns>" +
            "<w:pPr>" +
              "<w:pStyle w:val=`"ListParagraph`"/>" +
              "<w:keepNext/>" +
              "<w:keepLines/>" +
              "<w:numPr><w:ilvl w:val=`"2`"/><w:numId w:val=`"1`"/></w:numPr>" +
              "<w:ind w:left=`"1440`"/>" +
              "<w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr>" +
            "</w:pPr>" +
            "<w:r>" +
              "<w:rPr><w:strike/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr>" +
              "<w:t>Email the students in STAT 311 AA and AB reminding them to bring their laptops to class</w:t>" +
            "</w:r>" +
            "<w:r>" +
              "<w:rPr><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/></w:rPr>" +
              "<w:t>.</w:t>" +
            "</w:r>" +
          "</w:p>"

$noneTarget.InsertXML($newXml)

# ---------------------------------------------------------------------------
# Step 2: remove the <w:bookmarkStart w:id="8".../> from the paragraph that
# used to carry it (the page-break paragraph right before the "Thursday,
# March 28, 2024" heading) - it now only lives at the end of "None." above.
# ---------------------------------------------------------------------------

$thuRange2 = $d.Content
$thuRange2.Find.Execute("Thursday, March 28, 2024") | Out-Null
$thuStart2 = $thuRange2.Start

$pageBreakPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.End -eq $thuStart2) {
        $pageBreakPara = $p
    }
}
if ($null -eq $pageBreakPara) {
    throw "Could not find the page-break paragraph before Thursday heading"
}

$brStart = $pageBreakPara.Range.Start
$brEnd = $pageBreakPara.Range.End - 1
$brTarget = $d.Range($brStart, $brEnd)

$brXml = "<w:p $wns><w:r><w:br w:type=`"page`"/></w:r></w:p>"
$brTarget.InsertXML($brXml)
